$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price (column D) and 1h volume-change (column E) values.
# Column D values are forced to Text format so Excel does not reinterpret
# numeric-looking strings (e.g. "324.21", "28.565.68") as numbers/dates.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.565.68"
$ws.Range("E2").Value = "  +2.01%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.866.53"
$ws.Range("E3").Value = "  +1.98%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.21"
$ws.Range("E5").Value = "  +0.03%  "
$ws.Range("E6").Value = "  +0.15%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4604"
$ws.Range("E7").Value = "  -1.25%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3871"
$ws.Range("E8").Value = "  +0.18%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07853"
$ws.Range("E9").Value = "  -0.24%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9746"
$ws.Range("E10").Value = "  +1.61%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.90"
$ws.Range("E11").Value = "  +0.32%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.891.50"
$ws.Range("E12").Value = "  +4.24%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.968"
$ws.Range("E13").Value = "  +1.02%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.684"
$ws.Range("E14").Value = "  +0.28%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06942"
$ws.Range("E15").Value = "  +1.53%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.11"
$ws.Range("E16").Value = "  +0.95%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.003"
$ws.Range("E17").Value = "  +0.11%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009995"
$ws.Range("E18").Value = "  +0.78%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.76"
$ws.Range("E19").Value = "  +1.07%  "
$ws.Range("E20").Value = "  +0.06%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "28.566.70"
$ws.Range("E21").Value = "  +2.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.258"
$ws.Range("E22").Value = "  -1.04%  "
$ws.Range("E23").Value = "  +0.48%  "
$ws.Range("E24").Value = "  +1.00%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.105.08"
$ws.Range("E25").Value = "  +3.28%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "152.26"
$ws.Range("E26").Value = "  -1.07%  "
$ws.Range("E27").Value = "  +0.63%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.781"
$ws.Range("E28").Value = "  +1.22%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.982"
$ws.Range("E29").Value = "  +1.29%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "119.05"
$ws.Range("E30").Value = "  +1.15%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09324"
$ws.Range("E31").Value = "  +0.69%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9175"
$ws.Range("E32").Value = "  -1.78%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.269"
$ws.Range("E33").Value = "  -0.10%  "
$ws.Range("E34").Value = "  +1.37%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.322"
$ws.Range("E35").Value = "  +0.84%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05782"
$ws.Range("E36").Value = "  -1.36%  "
$ws.Range("E37").Value = "  +0.72%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02072"
$ws.Range("E38").Value = "  -3.19%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.711"
$ws.Range("E39").Value = "  -0.83%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5614"
$ws.Range("E40").Value = "  +0.66%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1780"
$ws.Range("E41").Value = "  +1.32%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.759"
$ws.Range("E42").Value = "  -0.88%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.07164"
$ws.Range("E43").Value = "  +2.05%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "11.66"
$ws.Range("E44").Value = "  +0.29%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5287"
$ws.Range("E45").Value = "  +0.67%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.149"
$ws.Range("E46").Value = "  +1.55%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.137"
$ws.Range("E47").Value = "  +3.21%  "
$ws.Range("E48").Value = "  +0.50%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "112.55"
$ws.Range("E49").Value = "  -0.33%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.415"
$ws.Range("E50").Value = "  +4.00%  "
$ws.Range("E51").Value = "  +0.15%  "
